$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2167.75
$ws.Range("I28").Value = 2167.75
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 2167.75
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -1682.75
$ws.Range("N28").Value = ""
$ws.Range("H32").Value = 1712.75
$ws.Range("J32").Value = 1916.6666
$ws.Range("L32").Value = 1916.6666
$ws.Range("N32").Value = -2568.6666
$ws.Range("H41").Value = 452.75
$ws.Range("I41").Value = 586
$ws.Range("J41").Value = 230.66667
$ws.Range("K41").Value = 586
$ws.Range("L41").Value = 230.66667
$ws.Range("M41").Value = -146
$ws.Range("N41").Value = -1110.66667
$ws.Range("H74").Value = 4562.5
$ws.Range("I74").Value = 4000
$ws.Range("K74").Value = 4000
$ws.Range("M74").Value = -3064
$ws.Range("H77").Value = 4562.5
$ws.Range("I77").Value = 4000
$ws.Range("K77").Value = 20000
$ws.Range("M77").Value = -15320
$ws.Range("H86").Value = 52375.375
$ws.Range("I86").Value = 1820.6
$ws.Range("J86").Value = 136633.33
$ws.Range("K86").Value = 1820.6
$ws.Range("L86").Value = 136633.33
$ws.Range("M86").Value = -697.5999999999999
$ws.Range("N86").Value = -138879.33
$ws.Range("H89").Value = 52375.375
$ws.Range("I89").Value = 1820.6
$ws.Range("J89").Value = 136633.33
$ws.Range("K89").Value = 9103
$ws.Range("L89").Value = 683166.6499999999
$ws.Range("M89").Value = -3487
$ws.Range("N89").Value = -694398.6499999999
$ws.Range("H132").Value = 2370.0356
$ws.Range("I132").Value = 2413.3704
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 7240.111199999999
$ws.Range("L132").Value = 3600
$ws.Range("M132").Value = -4710.111199999999
$ws.Range("N132").Value = -8660

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1500
$ws.Range("I45").Value = 1500
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1500
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1123
$ws.Range("N45").Value = ""
$ws.Range("H61").Value = 5986.615
$ws.Range("I61").Value = 4346.5454
$ws.Range("J61").Value = 15007
$ws.Range("K61").Value = 4346.5454
$ws.Range("L61").Value = 15007
$ws.Range("M61").Value = -4134.5454
$ws.Range("N61").Value = -15431
$ws.Range("H63").Value = 10392.083
$ws.Range("I63").Value = 12040.5
$ws.Range("J63").Value = 2150
$ws.Range("K63").Value = 12040.5
$ws.Range("L63").Value = 2150
$ws.Range("M63").Value = -11354.5
$ws.Range("N63").Value = -3522
$ws.Range("H66").Value = 10392.083
$ws.Range("I66").Value = 12040.5
$ws.Range("J66").Value = 2150
$ws.Range("K66").Value = 60202.5
$ws.Range("L66").Value = 10750
$ws.Range("M66").Value = -56770.5
$ws.Range("N66").Value = -17614
$ws.Range("H136").Value = 5986.615
$ws.Range("I136").Value = 4346.5454
$ws.Range("J136").Value = 15007
$ws.Range("K136").Value = 13039.6362
$ws.Range("L136").Value = 45021
$ws.Range("M136").Value = -10489.6362
$ws.Range("N136").Value = -50121

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 26333.334
$ws.Range("I82").Value = 3000
$ws.Range("J82").Value = 38000
$ws.Range("K82").Value = 3000
$ws.Range("L82").Value = 38000
$ws.Range("M82").Value = -2617
$ws.Range("N82").Value = -38766
$ws.Range("H85").Value = 26333.334
$ws.Range("I85").Value = 3000
$ws.Range("J85").Value = 38000
$ws.Range("K85").Value = 3000
$ws.Range("L85").Value = 38000
$ws.Range("M85").Value = -1674
$ws.Range("N85").Value = -40652
$ws.Range("H86").Value = 2780.5334
$ws.Range("I86").Value = 2884.6924
$ws.Range("J86").Value = 2103.5
$ws.Range("K86").Value = 2884.6924
$ws.Range("L86").Value = 2103.5
$ws.Range("M86").Value = -1761.6924
$ws.Range("N86").Value = -4349.5
$ws.Range("H89").Value = 2780.5334
$ws.Range("I89").Value = 2884.6924
$ws.Range("J89").Value = 2103.5
$ws.Range("K89").Value = 14423.462
$ws.Range("L89").Value = 10517.5
$ws.Range("M89").Value = -8807.462
$ws.Range("N89").Value = -21749.5
$ws.Range("H107").Value = 1401.3889
$ws.Range("I107").Value = 1062.8462
$ws.Range("J107").Value = 2281.6
$ws.Range("K107").Value = 1062.8462
$ws.Range("L107").Value = 2281.6
$ws.Range("M107").Value = 857.1538
$ws.Range("N107").Value = -6121.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 52926.668
$ws.Range("J52").Value = 52926.668
$ws.Range("L52").Value = 52926.668
$ws.Range("N52").Value = -53514.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1007.0323
$ws.Range("I5").Value = 593.9048
$ws.Range("K5").Value = 1781.7144
$ws.Range("M5").Value = -1669.7144
$ws.Range("H134").Value = 1308.9286
$ws.Range("I134").Value = 803.125
$ws.Range("J134").Value = 1983.3334
$ws.Range("K134").Value = 2409.375
$ws.Range("L134").Value = 5950.0002
$ws.Range("M134").Value = 2660.625
$ws.Range("N134").Value = -16090.0002
$ws.Range("H135").Value = 1007.0323
$ws.Range("I135").Value = 593.9048
$ws.Range("K135").Value = 5345.1432
$ws.Range("M135").Value = -2810.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 47500
$ws.Range("J18").Value = 47500
$ws.Range("L18").Value = 47500
$ws.Range("N18").Value = -48086
$ws.Range("H80").Value = 3500
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3500
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 3500
$ws.Range("M80").Value = ""
$ws.Range("N80").Value = -5496
$ws.Range("H83").Value = 3500
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3500
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 17500
$ws.Range("M83").Value = ""
$ws.Range("N83").Value = -27484
$ws.Range("H107").Value = 5251.25
$ws.Range("I107").Value = 6834
$ws.Range("J107").Value = 503
$ws.Range("K107").Value = 6834
$ws.Range("L107").Value = 503
$ws.Range("M107").Value = -4914
$ws.Range("N107").Value = -4343
$ws.Range("H113").Value = 1619.5
$ws.Range("I113").Value = 1619.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1619.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 550.5
$ws.Range("N113").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 8750
$ws.Range("J58").Value = 8750
$ws.Range("L58").Value = 8750
$ws.Range("N58").Value = -9270
$ws.Range("H61").Value = 2468.75
$ws.Range("I61").Value = 1867.5
$ws.Range("J61").Value = 3070
$ws.Range("K61").Value = 1867.5
$ws.Range("L61").Value = 3070
$ws.Range("M61").Value = -1665.5
$ws.Range("N61").Value = -3474
$ws.Range("H82").Value = 3134.3333
$ws.Range("I82").Value = 3900
$ws.Range("J82").Value = 2751.5
$ws.Range("K82").Value = 3900
$ws.Range("L82").Value = 2751.5
$ws.Range("M82").Value = -3539
$ws.Range("N82").Value = -3473.5
$ws.Range("H85").Value = 3134.3333
$ws.Range("I85").Value = 3900
$ws.Range("J85").Value = 2751.5
$ws.Range("K85").Value = 3900
$ws.Range("L85").Value = 2751.5
$ws.Range("M85").Value = -2652
$ws.Range("N85").Value = -5247.5
$ws.Range("H100").Value = 1628.875
$ws.Range("I100").Value = 1338.5
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 1338.5
$ws.Range("L100").Value = 2500
$ws.Range("M100").Value = -797.5
$ws.Range("N100").Value = -3582
$ws.Range("H113").Value = 2468.75
$ws.Range("I113").Value = 1867.5
$ws.Range("J113").Value = 3070
$ws.Range("K113").Value = 1867.5
$ws.Range("L113").Value = 3070
$ws.Range("M113").Value = 302.5
$ws.Range("N113").Value = -7410
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").Value = ""
$ws.Range("H133").Value = 21559.25
$ws.Range("J133").Value = 21559.25
$ws.Range("L133").Value = 21559.25
$ws.Range("N133").Value = -26619.25
$ws.Range("H136").Value = 5035.875
$ws.Range("I136").Value = 2696.4
$ws.Range("J136").Value = 8935
$ws.Range("K136").Value = 8089.200000000001
$ws.Range("L136").Value = 26805
$ws.Range("M136").Value = -5539.200000000001
$ws.Range("N136").Value = -31905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4000
$ws.Range("I62").Value = 3500
$ws.Range("K62").Value = 3500
$ws.Range("M62").Value = -2876
$ws.Range("H65").Value = 4000
$ws.Range("I65").Value = 3500
$ws.Range("K65").Value = 17500
$ws.Range("M65").Value = -14380
